$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newJson = '{"version": "1.2.0", "type": "plot", "attrs": {"style": "heatmap", "delimiter": "tab"}}'

$ws.Range("G14").Value = $newJson
$ws.Range("G16").Value = $newJson

$ws.Range("G16").Select()
